$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2:B52 (lamda_1) and C2:C52 (lamda_2) to new constant values for all data rows
$ws.Range("B2:B52").Value = 33.94444444444444
$ws.Range("C2:C52").Value = 1.95

# Update dic_nbre_clients_prob_poisson_2_values (column E) for changed rows
$ws.Range("E4").Value = 0.002
$ws.Range("E5").Value = 0.012
$ws.Range("E6").Value = 0.019
$ws.Range("E7").Value = 0.041
$ws.Range("E8").Value = 0.032
$ws.Range("E9").Value = 0.046
$ws.Range("E10").Value = 0.038
$ws.Range("E11").Value = 0.035
$ws.Range("E12").Value = 0.026
$ws.Range("E13").Value = 0.036
$ws.Range("E14").Value = 0.021
$ws.Range("E15").Value = 0.029
$ws.Range("E17").Value = 0.038
$ws.Range("E18").Value = 0.035
$ws.Range("E19").Value = 0.035
$ws.Range("E20").Value = 0.029
$ws.Range("E21").Value = 0.029
$ws.Range("E22").Value = 0.013
$ws.Range("E24").Value = 0.021
$ws.Range("E25").Value = 0.021
$ws.Range("E26").Value = 0.022
$ws.Range("E27").Value = 0.017
$ws.Range("E28").Value = 0.022
$ws.Range("E29").Value = 0.013
$ws.Range("E30").Value = 0.012
$ws.Range("E31").Value = 0.01
$ws.Range("E33").Value = 0.016
$ws.Range("E34").Value = 0.013
$ws.Range("E36").Value = 0.003
$ws.Range("E37").Value = 0.01
$ws.Range("E38").Value = 0.008
$ws.Range("E39").Value = 0.007
$ws.Range("E40").Value = 0.013
$ws.Range("E41").Value = 0.006
$ws.Range("E43").Value = 0.003
$ws.Range("E44").Value = 0.003
$ws.Range("E45").Value = 0.004
$ws.Range("E46").Value = 0.001
$ws.Range("E47").Value = 0.002
$ws.Range("E51").Value = 0.002

# Update dic_nbre_clients_poisson_2_keys (column D) for changed rows
$ws.Range("D46").Value = 46
$ws.Range("D47").Value = 48
$ws.Range("D48").Value = 49
$ws.Range("D49").Value = 50
$ws.Range("D50").Value = 54
$ws.Range("D51").Value = 57
$ws.Range("D52").Value = 75

# Remove rows that no longer exist in the updated distribution (53-58)
$ws.Range("A53:E58").EntireRow.Delete()
